# Update cryptos list (values + row 22/23/42/43 coin swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "36.582.35"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -1.61%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.065.21"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +0.43%  "
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  -0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "241.92"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -2.59%  "
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.660"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -0.94%  "
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +0.02%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "52.73"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -8.22%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "58.91"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -1.66%  "
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -6.96%  "
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -4.57%  "
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -0.73%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.899"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -1.23%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "14.60"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -9.39%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.362.39"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +0.18%  "
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -5.82%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.076.47"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +0.89%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "36.480.17"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -1.85%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "16.31"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -15.49%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "71.72"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -4.48%  "
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -4.14%  "
$c.Style = "Normal"

$c = $ws.Range("B22")
$c.NumberFormat = "@"
$c.Value = "Uniswap"
$c.Style = "Normal"

$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -3.78%  "
$c.Style = "Normal"

$c = $ws.Range("B23")
$c.NumberFormat = "@"
$c.Value = "BitcoinCash"
$c.Style = "Normal"

$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "236.54"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -0.62%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -0.12%  "
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -5.39%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.40"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  -2.25%  "
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -3.13%  "
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "163.93"
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -3.97%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "20.49"
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +1.17%  "
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -2.76%  "
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -1.90%  "
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -3.65%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.58"
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -1.37%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0596"
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -4.60%  "
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +0.02%  "
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -0.84%  "
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +2.90%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0825"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -6.44%  "
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -6.68%  "
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -5.11%  "
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -5.88%  "
$c.Style = "Normal"

$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = "ARBITRUM"
$c.Style = "Normal"

$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.13"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -3.10%  "
$c.Style = "Normal"

$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = "VeChain"
$c.Style = "Normal"

$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0215"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -3.79%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0934"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -6.49%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "94.13"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -2.69%  "
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.398.54"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +9.45%  "
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +8.32%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "15.56"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  -11.58%  "
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.34"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -3.79%  "
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.85"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -0.20%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.252.14"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +0.30%  "
$c.Style = "Normal"

